$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.03947
$ws.Range("H2").Value = 0.11841
$ws.Range("I2").Value = 0.02464674251283125
$ws.Range("J2").Value = 0.02464674251283126
$ws.Range("Q2").Value = 0.001065440023333333
$ws.Range("R2").Value = 0.009588960210000001
$ws.Range("S2").Value = 0.02464674251283125
$ws.Range("T2").Value = 0.02464674251283126

# Row 3
$ws.Range("I3").Value = 0.5082605406922069
$ws.Range("J3").Value = 0.5082605406922069
$ws.Range("S3").Value = 0.5082605406922069
$ws.Range("T3").Value = 0.5082605406922069

# Row 4
$ws.Range("G4").Value = 0.08494933333333334
$ws.Range("H4").Value = 0.254848
$ws.Range("I4").Value = 0.05304596770467038
$ws.Range("J4").Value = 0.05304596770467038
$ws.Range("Q4").Value = 0.002293093987555556
$ws.Range("R4").Value = 0.020637845888
$ws.Range("S4").Value = 0.05304596770467038
$ws.Range("T4").Value = 0.05304596770467038

# Row 5
$ws.Range("G5").Value = 0.434152
$ws.Range("H5").Value = 1.302456
$ws.Range("I5").Value = 0.2711029276774947
$ws.Range("J5").Value = 0.2711029276774947
$ws.Range("Q5").Value = 0.01171935437066667
$ws.Range("R5").Value = 0.105474189336
$ws.Range("S5").Value = 0.2711029276774947
$ws.Range("T5").Value = 0.2711029276774947

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.2289143333333333
$ws.Range("H6").Value = 0.686743
$ws.Range("I6").Value = 0.1429438214127968
$ws.Range("J6").Value = 0.1429438214127969
$ws.Range("Q6").Value = 0.006179237209222222
$ws.Range("R6").Value = 0.055613134883
$ws.Range("S6").Value = 0.1429438214127968
$ws.Range("T6").Value = 0.1429438214127969
